# Refresh the scraped cryptocurrency price/volume snapshot (GitHub Actions bot run).
# Column D = Price (text, as scraped), Column E = Volume(1h) change (text, "  +x.xx%  ").
# A leading apostrophe is used for numeric-looking Price values so Excel keeps
# them as literal text (preserving exact digits such as trailing/leading zeros)
# instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.189.02"
$ws.Range("E2").Value = "  +3.49%  "

$ws.Range("D3").Value = "2.062.14"
$ws.Range("E3").Value = "  +2.91%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E5").Value = "  +2.90%  "

$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  +1.91%  "

$ws.Range("D7").Value = "'58.56"
$ws.Range("E7").Value = "  +7.34%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +3.49%  "

$ws.Range("D10").Value = "'0.0810"
$ws.Range("E10").Value = "  +3.86%  "

$ws.Range("E11").Value = "  +0.90%  "

$ws.Range("D12").Value = "2.367.85"
$ws.Range("E12").Value = "  +3.02%  "

$ws.Range("D13").Value = "'14.68"
$ws.Range("E13").Value = "  +4.56%  "

$ws.Range("D14").Value = "'20.73"
$ws.Range("E14").Value = "  +3.29%  "

$ws.Range("D15").Value = "'0.754"
$ws.Range("E15").Value = "  +2.60%  "

$ws.Range("D16").Value = "'5.29"
$ws.Range("E16").Value = "  +4.07%  "

$ws.Range("D17").Value = "2.061.72"
$ws.Range("E17").Value = "  +3.31%  "

$ws.Range("D18").Value = "38.050.62"
$ws.Range("E18").Value = "  +3.23%  "

$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  +0.98%  "

$ws.Range("D20").Value = "'70.00"
$ws.Range("E20").Value = "  +2.28%  "

$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = "  +2.54%  "

$ws.Range("D22").Value = "'225.00"
$ws.Range("E22").Value = "  +1.26%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("E25").Value = "  +4.15%  "

$ws.Range("D26").Value = "'9.33"
$ws.Range("E26").Value = "  +2.52%  "

$ws.Range("D27").Value = "'166.36"
$ws.Range("E27").Value = "  +0.70%  "

$ws.Range("E28").Value = "  +7.78%  "

$ws.Range("D29").Value = "'19.09"
$ws.Range("E29").Value = "  +2.54%  "

$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("E31").Value = "  +2.51%  "

$ws.Range("D32").Value = "'4.56"
$ws.Range("E32").Value = "  +1.77%  "

$ws.Range("D33").Value = "'4.62"
$ws.Range("E33").Value = "  +5.41%  "

$ws.Range("D34").Value = "'0.0615"
$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("E35").Value = "  +8.01%  "

$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  +0.85%  "

$ws.Range("D37").Value = "'6.03"
$ws.Range("E37").Value = "  +15.58%  "

$ws.Range("E38").Value = "  +6.57%  "

$ws.Range("E39").Value = "  -0.10%  "

$ws.Range("E40").Value = "  +2.76%  "

$ws.Range("D41").Value = "'98.56"
$ws.Range("E41").Value = "  +4.36%  "

$ws.Range("D42").Value = "1.482.68"
$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("D43").Value = "'0.0947"
$ws.Range("E43").Value = "  +3.82%  "

$ws.Range("D44").Value = "'16.89"
$ws.Range("E44").Value = "  +4.89%  "

$ws.Range("E45").Value = "  +3.67%  "

$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("D47").Value = "'4.08"
$ws.Range("E47").Value = "  +17.58%  "

$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").Value = "'2.96"
$ws.Range("E49").Value = "  +2.06%  "

$ws.Range("E50").Value = "  -0.41%  "

$ws.Range("D51").Value = "2.255.04"
$ws.Range("E51").Value = "  +3.39%  "
